$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Remove the empty "ListParagraph" paragraph (spacing after=0,
#    ind left=1440) that sits just before the
#    "Create module to output memory..." bullet.
# ---------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text
    if ($txt -eq "Create module to output memory for options 1 and 3 of the menu.`r") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not locate target paragraph"
}

$prev = $target.Previous()
if ($prev.Range.Text -eq "`r") {
    $prev.Range.Delete()
}

# Re-find the target paragraph (collection indices shift after the delete).
$target = $null
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text
    if ($txt -eq "Create module to output memory for options 1 and 3 of the menu.`r") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not re-locate target paragraph after delete"
}

# ---------------------------------------------------------------
# 2) Replace that paragraph with the new block of five paragraphs:
#    - reworded "ddr2 RAM" bullet (with a lastRenderedPageBreak)
#    - blank spacer
#    - new microphone/speaker bullet (mic's/gramStart/gramEnd runs)
#    - blank spacer
#    - new submenu bullet split across three runs, ending with the
#      original "module to output memory..." sentence.
# ---------------------------------------------------------------
$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>Create a memory interfacing module to use the already written modules given for interfacing the ddr2 RAM</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:spacing w:after="0"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Create a microphone and speaker interfacing module to use the already written modules given for interfacing between </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>mic\u2019s</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> and speakers</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:spacing w:after="0"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Create </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">a submenu </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>module to output memory for options 1 and 3 of the menu.</w:t></w:r></w:p>
'@
$newXml = $newXml.Replace("\u2019", [char]0x2019)

$target.Range.InsertXML($newXml)

# ---------------------------------------------------------------
# 3) Add a lastRenderedPageBreak before "LGSB".
# ---------------------------------------------------------------
$lgsb = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "LGSB`r") {
        $lgsb = $p
        break
    }
}
if ($lgsb -eq $null) {
    throw "Could not locate LGSB paragraph"
}

$lgsbXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="144"/><w:szCs w:val="144"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="144"/><w:szCs w:val="144"/></w:rPr><w:lastRenderedPageBreak/><w:t>LGSB</w:t></w:r></w:p>'
$lgsb.Range.InsertXML($lgsbXml)

Write-Output "OK"
